$d = $word.ActiveDocument

$bm = $d.Bookmarks("_GoBack")
$goBackStart = $bm.Start
$scope = $d.Range($goBackStart - 120, $goBackStart + 20)

$hit = $d.Range($scope.Start, $scope.End)
$found = $hit.Find.Execute("TokenizeLine", $false, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
$tokenizeLineStart = $hit.Start
$tokenizeLineEnd = $hit.End
$tokenizeEnd = $tokenizeLineEnd - 4   # "Line" is the last 4 chars

$lineRange = $d.Range($tokenizeEnd, $tokenizeLineEnd)
$lineRange.Delete()

$insPoint = $d.Range($tokenizeEnd, $tokenizeEnd)
$insPoint.InsertAfter("String")

$d.Bookmarks.Add("__seam__", $d.Range($tokenizeEnd, $tokenizeEnd))
$d.Bookmarks("__seam__").Delete()

# try bracket split at very start of "Tokenize" run
$d.Bookmarks.Add("__edge1__", $d.Range($tokenizeLineStart, $tokenizeLineStart))
$d.Bookmarks("__edge1__").Delete()
